$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "EB000013"
$ws.Range("C3").Value = "EB000016"

$ws.Range("C7").Select() | Out-Null
